# Auto-generated PowerShell COM-interop script
# Appends newly observed daily snapshot rows to each of the three portfolio sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: append 11 rows ---
$data1 = @(
    @(127, "大智 (稳健智远)", "000333", "美的集团", 3.07, 42.89719527444501, 72.8, 3122.915815979597, 101775.9410260618, "202506241600"),
    @(128, "大智 (稳健智远)", "510050", "上证50ETF", 5.08, 1852.638869852698, 2.79, 5168.862446889028, 101775.9410260618, "202506241600"),
    @(129, "大智 (稳健智远)", "510300", "沪深300ETF", 4.94, 1276.881426590205, 3.94, 5030.912820765408, 101775.9410260618, "202506241600"),
    @(130, "大智 (稳健智远)", "518880", "黄金ETF", 4.88, 673.9096418114974, 7.37, 4966.714060150736, 101775.9410260618, "202506241600"),
    @(131, "大智 (稳健智远)", "600085", "同仁堂", 2, 56.06334956913254, 36.37, 2039.02402382935, 101775.9410260618, "202506241600"),
    @(132, "大智 (稳健智远)", "600900", "长江电力", 20.1, 669.0422707938175, 30.58, 20459.31264087494, 101775.9410260618, "202506241600"),
    @(133, "大智 (稳健智远)", "600989", "宝丰能源", 4.95, 308.5861230826723, 16.33, 5039.211389940038, 101775.9410260618, "202506241600"),
    @(134, "大智 (稳健智远)", "601899", "XD紫金矿", 9.859999999999999, 541.4194359293219, 18.54, 10037.91634212963, 101775.9410260618, "202506241600"),
    @(135, "大智 (稳健智远)", "HK02899", "紫金矿业", 9.869999999999999, 521.4694874201556, 19.26, 10043.5023277122, 101775.9410260618, "202506241600"),
    @(136, "大智 (稳健智远)", "HK06881", "中国银河", 5.21, 600.7968033130801, 8.82, 5299.027805221366, 101775.9410260618, "202506241600"),
    @(137, "大智 (稳健智远)", "100000", "现金", 30.04, 30568.54135256952, 1, 30568.54135256952, 101775.9410260618, "202506241600")
)
$ws1 = $wb.Worksheets.Item(1)
foreach ($row in $data1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).NumberFormat = "@"
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
    $ws1.Cells.Item($r, 6).Value = $row[6]
    $ws1.Cells.Item($r, 7).Value = $row[7]
    $ws1.Cells.Item($r, 8).Value = $row[8]
    $ws1.Cells.Item($r, 9).NumberFormat = "@"
    $ws1.Cells.Item($r, 9).Value = $row[9]
}

# --- Sheet 2: append 7 rows ---
$data2 = @(
    @(86, "大成 (锐进先锋)", "000725", "京东方A", 5.08, 1264.047242468624, 3.98, 5030.908025025124, 98961.38396064726, "202506241600"),
    @(87, "大成 (锐进先锋)", "159781", "科创创业ETF", 5.06, 9277.629760760277, 0.54, 5009.92007081055, 98961.38396064726, "202506241600"),
    @(88, "大成 (锐进先锋)", "513100", "纳指ETF", 5.03, 3131.938709046463, 1.59, 4979.782547383877, 98961.38396064726, "202506241600"),
    @(89, "大成 (锐进先锋)", "513290", "纳指生物科技ETF", 0.98, 870.2909333102562, 1.11, 966.0229359743844, 98961.38396064726, "202506241600"),
    @(90, "大成 (锐进先锋)", "603119", "浙江荣泰", 46.39, 1069.978093782073, 42.91, 45912.76000418875, 98961.38396064726, "202506241600"),
    @(91, "大成 (锐进先锋)", "688290", "景业智能", 7.64, 147.121987895035, 51.38, 7559.127738046899, 98961.38396064726, "202506241600"),
    @(92, "大成 (锐进先锋)", "100000", "现金", 29.81, 29502.86263921768, 1, 29502.86263921768, 98961.38396064726, "202506241600")
)
$ws2 = $wb.Worksheets.Item(2)
foreach ($row in $data2) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).NumberFormat = "@"
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 4).Value = $row[4]
    $ws2.Cells.Item($r, 5).Value = $row[5]
    $ws2.Cells.Item($r, 6).Value = $row[6]
    $ws2.Cells.Item($r, 7).Value = $row[7]
    $ws2.Cells.Item($r, 8).Value = $row[8]
    $ws2.Cells.Item($r, 9).NumberFormat = "@"
    $ws2.Cells.Item($r, 9).Value = $row[9]
}

# --- Sheet 3: append 15 rows ---
$data3 = @(
    @(181, "范式进化投资组合", "000333", "美的集团", 1.02, 14.02515882310653, 72.8, 1021.031562322155, 100033.6101130758, "202506241600"),
    @(182, "范式进化投资组合", "000725", "京东方A", 5.11, 1284.618017653691, 3.98, 5112.77971026169, 100033.6101130758, "202506241600"),
    @(183, "范式进化投资组合", "159781", "科创创业ETF", 5.09, 9428.611488061997, 0.54, 5091.450203553479, 100033.6101130758, "202506241600"),
    @(184, "范式进化投资组合", "510050", "上证50ETF", 5.07, 1817.150577699221, 2.79, 5069.850111780826, 100033.6101130758, "202506241600"),
    @(185, "范式进化投资组合", "510300", "沪深300ETF", 4.93, 1252.422077361618, 3.94, 4934.542984804775, 100033.6101130758, "202506241600"),
    @(186, "范式进化投资组合", "513100", "纳指ETF", 1.01, 636.5814125697908, 1.59, 1012.164445985967, 100033.6101130758, "202506241600"),
    @(187, "范式进化投资组合", "513290", "纳指生物科技ETF", 0.98, 884.4538210040457, 1.11, 981.7437413144909, 100033.6101130758, "202506241600"),
    @(188, "范式进化投资组合", "518880", "黄金ETF", 0.97, 132.2001081659486, 7.37, 974.3147971830412, 100033.6101130758, "202506241600"),
    @(189, "范式进化投资组合", "600085", "同仁堂", 1, 27.49471300507762, 36.37, 999.982711994673, 100033.6101130758, "202506241600"),
    @(190, "范式进化投资组合", "600900", "长江电力", 1, 32.81132034584936, 30.58, 1003.370176176073, 100033.6101130758, "202506241600"),
    @(191, "范式进化投资组合", "600989", "宝丰能源", 4.94, 302.6749902285196, 16.33, 4942.682590431725, 100033.6101130758, "202506241600"),
    @(192, "范式进化投资组合", "601899", "XD紫金矿", 9.84, 531.0482559694855, 18.54, 9845.634665674261, 100033.6101130758, "202506241600"),
    @(193, "范式进化投资组合", "603119", "浙江荣泰", 1.04, 24.16423640557475, 42.91, 1036.887384163213, 100033.6101130758, "202506241600"),
    @(194, "范式进化投资组合", "HK06881", "中国银河", 1.04, 117.8576436007749, 8.82, 1039.504416558835, 100033.6101130758, "202506241600"),
    @(195, "范式进化投资组合", "100000", "现金", 56.95, 56967.67061087058, 1, 56967.67061087058, 100033.6101130758, "202506241600")
)
$ws3 = $wb.Worksheets.Item(3)
foreach ($row in $data3) {
    $r = $row[0]
    $ws3.Cells.Item($r, 1).Value = $row[1]
    $ws3.Cells.Item($r, 2).NumberFormat = "@"
    $ws3.Cells.Item($r, 2).Value = $row[2]
    $ws3.Cells.Item($r, 3).Value = $row[3]
    $ws3.Cells.Item($r, 4).Value = $row[4]
    $ws3.Cells.Item($r, 5).Value = $row[5]
    $ws3.Cells.Item($r, 6).Value = $row[6]
    $ws3.Cells.Item($r, 7).Value = $row[7]
    $ws3.Cells.Item($r, 8).Value = $row[8]
    $ws3.Cells.Item($r, 9).NumberFormat = "@"
    $ws3.Cells.Item($r, 9).Value = $row[9]
}
